$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "ODI Bowling Extra" worksheet at the end of the workbook
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $battingExtra)
$ws5.Name = "ODI Bowling Extra"

# Match the page margins used by the rest of the workbook
$ps = $ws5.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2) Populate the new sheet with MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL
# ---------------------------------------------------------------------------
$dataRange = $ws5.Range("A1:C21")
$dataRange.NumberFormat = "@"

$ws5.Cells.Item(1, 1).Value = "MATCH_CODE"
$ws5.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$ws5.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"

$ws5.Cells.Item(2, 1).Value = "4126"
$ws5.Cells.Item(2, 2).Value = ""
$ws5.Cells.Item(2, 3).Value = ""
$ws5.Cells.Item(3, 1).Value = "4226"
$ws5.Cells.Item(3, 2).Value = "0"
$ws5.Cells.Item(3, 3).Value = ""
$ws5.Cells.Item(4, 1).Value = "4303"
$ws5.Cells.Item(4, 2).Value = ""
$ws5.Cells.Item(4, 3).Value = ""
$ws5.Cells.Item(5, 1).Value = "4307"
$ws5.Cells.Item(5, 2).Value = ""
$ws5.Cells.Item(5, 3).Value = ""
$ws5.Cells.Item(6, 1).Value = "4334"
$ws5.Cells.Item(6, 2).Value = "0"
$ws5.Cells.Item(6, 3).Value = "10.00%"
$ws5.Cells.Item(7, 1).Value = "4459"
$ws5.Cells.Item(7, 2).Value = "0"
$ws5.Cells.Item(7, 3).Value = ""
$ws5.Cells.Item(8, 1).Value = "4460"
$ws5.Cells.Item(8, 2).Value = "0"
$ws5.Cells.Item(8, 3).Value = "20.00%"
$ws5.Cells.Item(9, 1).Value = "4487"
$ws5.Cells.Item(9, 2).Value = "0"
$ws5.Cells.Item(9, 3).Value = "10.00%"
$ws5.Cells.Item(10, 1).Value = "4491"
$ws5.Cells.Item(10, 2).Value = "0"
$ws5.Cells.Item(10, 3).Value = "10.00%"
$ws5.Cells.Item(11, 1).Value = "4524"
$ws5.Cells.Item(11, 2).Value = ""
$ws5.Cells.Item(11, 3).Value = ""
$ws5.Cells.Item(12, 1).Value = "4526"
$ws5.Cells.Item(12, 2).Value = "0"
$ws5.Cells.Item(12, 3).Value = "10.00%"
$ws5.Cells.Item(13, 1).Value = "4529"
$ws5.Cells.Item(13, 2).Value = "0"
$ws5.Cells.Item(13, 3).Value = ""
$ws5.Cells.Item(14, 1).Value = "4619"
$ws5.Cells.Item(14, 2).Value = ""
$ws5.Cells.Item(14, 3).Value = ""
$ws5.Cells.Item(15, 1).Value = "4657"
$ws5.Cells.Item(15, 2).Value = "0"
$ws5.Cells.Item(15, 3).Value = ""
$ws5.Cells.Item(16, 1).Value = "4698"
$ws5.Cells.Item(16, 2).Value = "0"
$ws5.Cells.Item(16, 3).Value = ""
$ws5.Cells.Item(17, 1).Value = "4699"
$ws5.Cells.Item(17, 2).Value = "0"
$ws5.Cells.Item(17, 3).Value = "10.00%"
$ws5.Cells.Item(18, 1).Value = "4700"
$ws5.Cells.Item(18, 2).Value = "0"
$ws5.Cells.Item(18, 3).Value = ""
$ws5.Cells.Item(19, 1).Value = "4731"
$ws5.Cells.Item(19, 2).Value = "0"
$ws5.Cells.Item(19, 3).Value = "10.00%"
$ws5.Cells.Item(20, 1).Value = "4746"
$ws5.Cells.Item(20, 2).Value = "0"
$ws5.Cells.Item(20, 3).Value = "10.00%"
$ws5.Cells.Item(21, 1).Value = "4751"
$ws5.Cells.Item(21, 2).Value = ""
$ws5.Cells.Item(21, 3).Value = ""

# Drop the temporary "text" number format again so the data rows end up
# without any explicit style, same as the other generated sheets
$dataRange.Style = "Normal"

# Re-apply the bold / bordered header look used by the other sheets in this
# workbook by copying the existing header format onto the new header row
$headerSource = $battingExtra.Range("A1:C1")
$headerSource.Copy()
$ws5.Range("A1:C1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) "ODI Batting Extra": drop the attribute cells that never actually held
#    any scraped value so only the cells with real data remain
# ---------------------------------------------------------------------------
$battingExtra.Range("B6:E6").ClearContents()
$battingExtra.Range("E9").ClearContents()
$battingExtra.Range("B10:E10").ClearContents()
$battingExtra.Range("E11").ClearContents()
$battingExtra.Range("B12:E12").ClearContents()
$battingExtra.Range("B13:E13").ClearContents()
$battingExtra.Range("B21:E21").ClearContents()

# Restore the original active sheet
$wb.Worksheets.Item(1).Activate()
